$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.758.76"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.803.13"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.91"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.40"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000248"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.87"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.445.26"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.832.86"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.46"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.779.03"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.06"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.51"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.79"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000146"
$ws.Range("E23").Value = "  -4.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.06"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.05"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.952.90"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.27"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.04"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0995"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.77"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  -3.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.21"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.71"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.298"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.77"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.87"
$ws.Range("E47").Value = "  +4.80%  "
$ws.Range("E48").Value = "  +11.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.34"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "395.12"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  +1.93%  "
